$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.490.67'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.39%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.553.31'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.79%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.34%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.87'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("E6").Value = '  -1.66%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.36%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.13'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.243'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.80%  '

# Row 10
$ws.Range("E10").Value = '  -0.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0893'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.773.51'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.91%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.552.80'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.92%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.451.91'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.14%  '

# Row 15
$ws.Range("E15").Value = '  -1.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.510'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.35%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.14'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.39%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.36'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.20%  '

# Row 20
$ws.Range("E20").Value = '  -2.35%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.23%  '

# Row 22
$ws.Range("E22").Value = '  -0.70%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.89'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.50%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.34%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.76'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.82%  '

# Row 27
$ws.Range("E27").Value = '  -1.28%  '

# Row 28
$ws.Range("E28").Value = '  -0.34%  '

# Row 29
$ws.Range("E29").Value = '  -2.84%  '

# Row 30
$ws.Range("E30").Value = '  -2.98%  '

# Row 31
$ws.Range("E31").Value = '  -4.52%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.16'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.387.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.88%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.00'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.03'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.60%  '

# Row 36
$ws.Range("E36").Value = '  -1.47%  '

# Row 37
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.67'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.33%  '

# Row 38
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.95%  '

# Row 39
$ws.Range("E39").Value = '  -1.58%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.95'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.12%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.513'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.45%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.27%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.771'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0460'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.35'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.65%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.73'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.687.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.86%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.870'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.87%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.38'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.21'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.57%  '

# Row 51
$ws.Range("E51").Value = '  -0.55%  '
